# Adds missing data in metadata: attribute_label ("steelhead") for every
# attribute row, plus units/min/max for the date attribute and the various
# steelhead/chinook/lamprey/sucker count attributes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("attribute")
$ws.Activate()

# ---- date attribute (row 2) ----
$ws.Range("D2").Value = "steelhead"
$ws.Range("J2").Value = "YYYY-MM-DD"

$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = "2021-02-02"
$ws.Range("M2").NumberFormat = "@"
$ws.Range("M2").Value = "2022-04-27"

# ---- section / riffle attributes (rows 3-4) ----
$ws.Range("D3").Value = "steelhead"
$ws.Range("D4").Value = "steelhead"

# ---- count attributes (rows 5-14): adult_sh, juve_sh, sh_redds, sh_carcass,
#      live_chn, chn_redd, live_pl, pl_redd, sasu_redd, unk_redd ----
$countRows = @{
    5  = 3
    6  = 8
    7  = 1
    8  = 1
    9  = 2
    10 = 2
    11 = 2
    12 = 2
    13 = 431
    14 = 1
}

foreach ($row in $countRows.Keys) {
    $ws.Range("D$row").Value = "steelhead"
    $ws.Range("G$row").Value = "count"
    $ws.Range("L$row").Value = 0
    $ws.Range("M$row").Value = $countRows[$row]
}

# ---- view niceties matching the saved workbook (zoom + active cell) ----
$excel.ActiveWindow.Zoom = 125
$ws.Range("D3").Select()
